$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date label for row 103 must be stored as text (shared string), not
# auto-converted to a date serial number. Force a text format, assign the
# value, then restore the default "Normal" style so no stray number format
# sticks to the cell (matches the plain, unstyled cells used elsewhere in
# this row-band).
$ws.Range("A103").NumberFormat = "@"
$ws.Range("A103").Value = "01-04-2021"
$ws.Range("A103").Style = "Normal"

# New row of data (row 103), columns B..Z
$newRow = @(134, 131.9, 131.1, 122.8, 135.3, 128.7, 137.1, 140.9, 142.8, 139.1, 165.6, 169.4, 131.4, 197.3, 199.2, 180.9, 138.2, 126.9, 127.3, 125.1, 117.7, 101.1, 129, 140.1, 144.7)

for ($i = 0; $i -lt $newRow.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item(103, $col).Value = $newRow[$i]
}

# Revised values in row 102
$row102updates = @{
    2  = 138.7
    5  = 121.3
    8  = 178.3
    9  = 138.8
    10 = 139.9
    11 = 138.7
    12 = 157.7
    13 = 160.1
    14 = 121.9
    15 = 195
    16 = 197.5
    17 = 171.2
    18 = 132.7
    19 = 135
    20 = 125.8
    21 = 123.5
    22 = 122.9
    23 = 94.40000000000001
    24 = 128
    25 = 141.1
    26 = 148.2
}

foreach ($col in $row102updates.Keys) {
    $ws.Cells.Item(102, $col).Value = $row102updates[$col]
}
